$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J12").Value = 0.012
$ws.Range("K12").Value = 0.109
$ws.Range("L12").Value = 11.96
$ws.Range("M12").Value = 6.788

$ws.Range("J13").Value = 0.002
$ws.Range("K13").Value = 0.185
$ws.Range("L13").Value = 12
$ws.Range("M13").Value = 6.861

$ws.Range("J14").Value = 54.73
$ws.Range("K14").Value = 0.195
$ws.Range("L14").Value = 12.06
$ws.Range("M14").Value = 8.38

$ws.Range("M19").Select()
